$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -19.2308319290313
$ws.Range("B2").Value = -6.97169899366997
$ws.Range("A3").Value = -23.9802852812849
$ws.Range("B3").Value = -7.7982463176201
$ws.Range("A4").Value = -16.960506082904
$ws.Range("B4").Value = -5.21503806213296
$ws.Range("A5").Value = -6.0557950345974
$ws.Range("B5").Value = 8.66878944083269
$ws.Range("A6").Value = -8.3173126196345
$ws.Range("B6").Value = 7.72200995182726
$ws.Range("A7").Value = -11.7390376770988
$ws.Range("B7").Value = 6.55014202270683
$ws.Range("A8").Value = -30.9038875155909
$ws.Range("B8").Value = -5.67362107160279
$ws.Range("A9").Value = -29.7553818137461
$ws.Range("B9").Value = -5.54597823818038
$ws.Range("A10").Value = -31.7424471672916
$ws.Range("B10").Value = -4.05570839629435
$ws.Range("A11").Value = 26.0426481467496
$ws.Range("B11").Value = -31.5358823004391
$ws.Range("A12").Value = 25.2837058303201
$ws.Range("B12").Value = -31.278645391384
$ws.Range("A13").Value = 25.4608257560563
$ws.Range("B13").Value = -29.8284757647492
$ws.Range("A14").Value = 14.0268080440438
$ws.Range("B14").Value = 24.1892809280537
$ws.Range("A15").Value = 15.4484645424469
$ws.Range("B15").Value = 24.7659992746942
$ws.Range("A16").Value = 14.7996743493779
$ws.Range("B16").Value = 25.3144793526506
$ws.Range("A17").Value = 19.2306529098194
$ws.Range("B17").Value = 10.0243896118943
$ws.Range("A18").Value = 18.5855206174471
$ws.Range("B18").Value = 9.92749469351463
$ws.Range("A19").Value = 19.8071849249185
$ws.Range("B19").Value = 10.7407092598987
